# Insumo "prima cedida" requiere fecha de expedicion para cruzar con gasto:
# add a new componente row for "camara_soat" (CAM) below the existing table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "camara_soat"
$ws.Range("B15").Value = "CAM"

# Leave the new row's clasificacion_adicional_cd (column C) selected/active,
# matching the cursor position left after the edit.
$ws.Range("C15").Select()
